$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Job")

# New result-status column (C) on the "Job" sheet: both submitted jobs come
# back as "Job Title Already exit" in this run.
$ws.Range("C1").Value = "Job Title Already exit"
$ws.Range("C2").Value = "Job Title Already exit"
